# Apply the 2023-01-04 20:49:24 re-crawl update to bread_coop_2023-01-04.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2023-01-04 12:56:28"
$newTimestamp = "2023-01-04 20:49:24"

# Refresh the timestamp (column O) for every data row (2..399)
for ($row = 2; $row -le 399; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}

# Row 11 (id 4066469 Naturaplan Bio Vollkorntoast 10 Scheiben): now out of stock online
$ws.Cells.Item(11, 13).Value = "Naturaplan Bio Vollkorntoast 10 Scheiben - Online kein Bestand 2.50 Schweizer Franken"

# Row 17: ratingAmount went from 21 to 22
$ws.Cells.Item(17, 4).Value = 22

# Row 33: ratingAmount went from 7 to 8
$ws.Cells.Item(33, 4).Value = 8

# Row 35: ratingAmount went from 41 to 42
$ws.Cells.Item(35, 4).Value = 42

# Row 53: ratingAmount went from 14 to 15
$ws.Cells.Item(53, 4).Value = 15

# Row 105 (Betty Bossi Bio Frischback Rustico Buttergipfel): now out of stock online
$ws.Cells.Item(105, 13).Value = "Betty Bossi Bio Frischback Rustico Buttergipfel - Online kein Bestand 3.50 Schweizer Franken"

# Row 194: ratingAmount went from 5 to 6
$ws.Cells.Item(194, 4).Value = 6

# Row 207: ratingAmount went from 66 to 67
$ws.Cells.Item(207, 4).Value = 67

# Row 398 (Aprikosentorte): ratingAmount was previously blank, ratingValue was 0 -> now 1 / 2
$ws.Cells.Item(398, 4).Value = 1
$ws.Cells.Item(398, 5).Value = 2
